$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header date string: "15-Apr-2025" -> "16-Apr-2025" (appears in A3 and B6) ---
$ws.Range("A3").Value = "1-Jul-2024 to 16-Apr-2025"
$ws.Range("B6").Value = "1-Jul-2024 to 16-Apr-2025"

# --- Row 17: item now has no Sales/Rate/Value entries; clear C:E and match the blank-row style ---
$ws.Range("B17").Value = 194
$ws.Range("C17:E17").ClearContents()
$ws.Range("C81:E81").Copy()
$ws.Range("C17:E17").PasteSpecial(-4122)

# --- Row 148: item now HAS Sales/Rate/Value entries; fill C:E and match the value-row style ---
$ws.Range("C149:E149").Copy()
$ws.Range("C148:E148").PasteSpecial(-4122)
$ws.Range("C148").Value = 0.5
$ws.Range("D148").Value = 18.440000000000001
$ws.Range("E148").Value = 9.2200000000000006

# --- Row 328: item now has no Sales/Rate/Value entries; clear C:E and match the blank-row style ---
$ws.Range("B328").Value = 79
$ws.Range("C328:E328").ClearContents()
$ws.Range("C81:E81").Copy()
$ws.Range("C328:E328").PasteSpecial(-4122)

# --- Remaining rows: straightforward quantity / rate / value updates (style unchanged) ---
$ws.Range("B9").Value = 36
$ws.Range("C9").Value = 480.5
$ws.Range("E9").Value = 480.5

$ws.Range("B11").Value = 176
$ws.Range("C11").Value = 193.5
$ws.Range("E11").Value = 175.68

$ws.Range("B18").Value = 87
$ws.Range("C18").Value = 9
$ws.Range("E18").Value = 18

$ws.Range("B20").Value = 122
$ws.Range("C20").Value = 21
$ws.Range("E20").Value = 42

$ws.Range("B23").Value = 169
$ws.Range("C23").Value = 109.75
$ws.Range("E23").Value = 219.5

$ws.Range("B36").Value = 23
$ws.Range("C36").Value = 35
$ws.Range("E36").Value = 91

$ws.Range("B37").Value = 81
$ws.Range("C37").Value = 23.5
$ws.Range("E37").Value = 61.1

$ws.Range("B39").Value = 104
$ws.Range("C39").Value = 10.5
$ws.Range("E39").Value = 28.35

$ws.Range("B46").Value = 44
$ws.Range("C46").Value = 17
$ws.Range("E46").Value = 39.95

$ws.Range("B49").Value = 416
$ws.Range("C49").Value = 333.5
$ws.Range("E49").Value = 373.95

$ws.Range("B50").Value = 318
$ws.Range("C50").Value = 478.5
$ws.Range("D50").Value = 1.12
$ws.Range("E50").Value = 538.31

$ws.Range("B52").Value = 296
$ws.Range("C52").Value = 529
$ws.Range("E52").Value = 595.06

$ws.Range("B54").Value = 387
$ws.Range("C54").Value = 744
$ws.Range("E54").Value = 836.4

$ws.Range("B59").Value = 82
$ws.Range("C59").Value = 92.5
$ws.Range("E59").Value = 175.75

$ws.Range("B60").Value = 62
$ws.Range("C60").Value = 246
$ws.Range("E60").Value = 467.4

$ws.Range("B61").Value = 93
$ws.Range("C61").Value = 193.5
$ws.Range("E61").Value = 406.35

$ws.Range("B69").Value = 149
$ws.Range("C69").Value = 76.5
$ws.Range("E69").Value = 53.55

$ws.Range("B72").Value = 330
$ws.Range("C72").Value = 1116
$ws.Range("E72").Value = 1729.8

$ws.Range("B78").Value = 163
$ws.Range("C78").Value = 159

$ws.Range("B80").Value = 285
$ws.Range("C80").Value = -89

$ws.Range("B82").Value = 211
$ws.Range("C82").Value = 104.5

$ws.Range("B83").Value = 250
$ws.Range("C83").Value = 33.5
$ws.Range("E83").Value = 38.53

$ws.Range("B85").Value = 364
$ws.Range("C85").Value = 2675
$ws.Range("E85").Value = 3745

$ws.Range("B86").Value = 210
$ws.Range("C86").Value = 2075
$ws.Range("E86").Value = 2905

$ws.Range("B91").Value = 196
$ws.Range("C91").Value = 156
$ws.Range("E91").Value = 192.84

$ws.Range("B130").Value = 4
$ws.Range("C130").Value = 16
$ws.Range("E130").Value = 248

$ws.Range("B156").Value = 240
$ws.Range("C156").Value = 38.35
$ws.Range("E156").Value = 130.39

$ws.Range("B191").Value = 152
$ws.Range("C191").Value = 34.5
$ws.Range("E191").Value = 135.14

$ws.Range("B203").Value = 41
$ws.Range("C203").Value = 35.5
$ws.Range("E203").Value = 119.66

$ws.Range("B209").Value = 151
$ws.Range("C209").Value = 13
$ws.Range("E209").Value = 49.4

$ws.Range("B213").Value = 135
$ws.Range("C213").Value = 24.5
$ws.Range("E213").Value = 100.73

$ws.Range("B230").Value = 45
$ws.Range("C230").Value = 5.5
$ws.Range("E230").Value = 28.88

$ws.Range("B240").Value = 22
$ws.Range("C240").Value = 16.5
$ws.Range("E240").Value = 111.38

$ws.Range("B257").Value = 31
$ws.Range("C257").Value = 14
$ws.Range("E257").Value = 58.52

$ws.Range("B339").Value = 19
$ws.Range("C339").Value = 10.5
$ws.Range("E339").Value = 113.06

$ws.Range("B434").Value = 40
$ws.Range("C434").Value = 0.4
$ws.Range("E434").Value = 1.22

$ws.Range("B486").Value = 68
$ws.Range("C486").Value = 4.5
$ws.Range("E486").Value = 42.75

$ws.Range("B493").Value = 119
$ws.Range("C493").Value = 170
$ws.Range("E493").Value = 187

$ws.Range("B499").Value = 3
$ws.Range("C499").Value = 20
$ws.Range("E499").Value = 61

$ws.Range("B509").Value = 72
$ws.Range("C509").Value = 38.5
$ws.Range("E509").Value = 65.84

$ws.Range("B511").Value = 217
$ws.Range("C511").Value = 26
$ws.Range("E511").Value = 59.28

$ws.Range("B598").Value = 96
$ws.Range("C598").Value = 62.5
$ws.Range("E598").Value = 120.03

$ws.Range("B599").Value = 99
$ws.Range("C599").Value = 360.5
$ws.Range("E599").Value = 688.34

$ws.Range("B606").Value = 96
$ws.Range("C606").Value = 352
$ws.Range("E606").Value = 492.02

$ws.Range("B627").Value = 207
$ws.Range("C627").Value = 126.5
$ws.Range("E627").Value = 88.55

$ws.Range("B628").Value = 157
$ws.Range("C628").Value = 65.5
$ws.Range("E628").Value = 44.9

$ws.Range("C638").Value = 34072.43
$ws.Range("E638").Value = 60253.16

$excel.CutCopyMode = $false
